# Apply year-over-year financial figure updates to the LMDCF sheet
# (commit: "Doing Updates for Financials")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LMDCF")

# Row 8
$ws.Range("E8").Value = 2400
$ws.Range("F8").Value = 3700
$ws.Range("J8").Value = 1500
# Row 9
$ws.Range("H9").Value = 100
# Row 10
$ws.Range("D10").Value = 1900
$ws.Range("E10").Value = 2100
$ws.Range("F10").Value = 3400
$ws.Range("H10").Value = 1300
$ws.Range("J10").Value = 1400
# Row 12
$ws.Range("D12").Value = 2000
# Row 14
$ws.Range("D14").Value = 1600
# Row 15
$ws.Range("F15").Value = 500
$ws.Range("G15").Value = 400
$ws.Range("J15").Value = 1900
# Row 17
$ws.Range("D17").Value = 6400
$ws.Range("F17").Value = 1700
$ws.Range("H17").Value = 1200
$ws.Range("I17").Value = 2200
$ws.Range("J17").Value = 4700
# Row 18
$ws.Range("D18").Value = -4300
$ws.Range("F18").Value = 1900
$ws.Range("I18").Value = -700
$ws.Range("J18").Value = -3100
# Row 21
$ws.Range("D21").Value = -3700
$ws.Range("F21").Value = 2800
$ws.Range("I21").Value = -400
# Row 22
$ws.Range("J22").Value = 200
# Row 23
$ws.Range("D23").Value = -4500
$ws.Range("F23").Value = 2100
$ws.Range("I23").Value = -800
$ws.Range("J23").Value = -3400
# Row 26
$ws.Range("D26").Value = -4700
$ws.Range("J26").Value = -3500
# Row 27
$ws.Range("D27").Value = -4700
$ws.Range("J27").Value = -3500
# Row 33
$ws.Range("D33").Value = -4700
$ws.Range("J33").Value = -3500
# Row 35
$ws.Range("D35").Value = -4700
$ws.Range("J35").Value = -3500
# Row 41
$ws.Range("D41").Value = 200
# Row 43
$ws.Range("G43").Value = 600
$ws.Range("H43").Value = 700
# Row 46
$ws.Range("D46").Value = 1100
$ws.Range("F46").Value = 2100
$ws.Range("J46").Value = 1300
# Row 49
$ws.Range("E49").Value = 2300
$ws.Range("F49").Value = 1700
$ws.Range("G49").Value = 700
$ws.Range("J49").Value = 900
# Row 54
$ws.Range("D54").Value = 1100
$ws.Range("E54").Value = 5300
$ws.Range("F54").Value = 3900
$ws.Range("G54").Value = 1800
$ws.Range("H54").Value = 1600
# Row 59
$ws.Range("H59").Value = 400
$ws.Range("J59").Value = 200
# Row 60
$ws.Range("D60").Value = 700
$ws.Range("E60").Value = 500
$ws.Range("G60").Value = 1200
# Row 66
$ws.Range("D66").Value = 700
$ws.Range("E66").Value = 500
$ws.Range("G66").Value = 1200
# Row 72
$ws.Range("D72").Value = -18500
$ws.Range("E72").Value = -13800
$ws.Range("F72").Value = -12800
$ws.Range("G72").Value = -14700
$ws.Range("H72").Value = -14800
$ws.Range("I72").Value = -14900
$ws.Range("J72").Value = -13800
# Row 76
$ws.Range("E76").Value = 4800
$ws.Range("F76").Value = 3000
# Row 81
$ws.Range("D81").Value = -4700
$ws.Range("J81").Value = -3500
# Row 83
$ws.Range("F83").Value = 500
$ws.Range("G83").Value = 400
$ws.Range("J83").Value = 1900
# Row 89
$ws.Range("F89").Value = 1100
# Row 91
$ws.Range("E91").Value = 0
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
